$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1773.6364
$ws.Range("J17").Value = 1990
$ws.Range("L17").Value = 5970
$ws.Range("N17").Value = -6306
$ws.Range("H40").Value = 1673.4286
$ws.Range("I40").Value = 1404.9166
$ws.Range("K40").Value = 1404.9166
$ws.Range("M40").Value = -1229.9166
$ws.Range("H43").Value = 500
$ws.Range("I43").Value = 500
$ws.Range("J43").Value = 500
$ws.Range("K43").Value = 500
$ws.Range("L43").Value = 500
$ws.Range("M43").Value = -431
$ws.Range("N43").Value = -638
$ws.Range("H138").Value = 3705176.2
$ws.Range("I138").Value = 1095.5897
$ws.Range("J138").Value = 13335786
$ws.Range("K138").Value = 3286.7691
$ws.Range("L138").Value = 40007358
$ws.Range("M138").Value = 1853.2309
$ws.Range("N138").Value = -40017638

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 107.333336
$ws.Range("I5").Value = 120
$ws.Range("J5").Value = 101
$ws.Range("K5").Value = 120
$ws.Range("L5").Value = 101
$ws.Range("M5").Value = -8
$ws.Range("N5").Value = -325
$ws.Range("H35").Value = 537
$ws.Range("I35").Value = 537
$ws.Range("K35").Value = 537
$ws.Range("M35").Value = -131
$ws.Range("H110").Value = 3430
$ws.Range("I110").Value = 3125
$ws.Range("K110").Value = 3125
$ws.Range("M110").Value = -1080
$ws.Range("H132").Value = 1746.6111
$ws.Range("I132").Value = 1266.3334
$ws.Range("J132").Value = 2707.1667
$ws.Range("K132").Value = 3799.0002
$ws.Range("L132").Value = 8121.500100000001
$ws.Range("M132").Value = -1269.0002
$ws.Range("N132").Value = -13181.5001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 107.333336
$ws.Range("I4").Value = 120
$ws.Range("J4").Value = 101
$ws.Range("K4").Value = 120
$ws.Range("L4").Value = 101
$ws.Range("M4").Value = -5
$ws.Range("N4").Value = -331
$ws.Range("H37").Value = 10459.333
$ws.Range("I37").Value = 3424.6667
$ws.Range("J37").Value = 13976.667
$ws.Range("K37").Value = 3424.6667
$ws.Range("L37").Value = 13976.667
$ws.Range("M37").Value = -3287.6667
$ws.Range("N37").Value = -14250.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 100.94118
$ws.Range("I7").Value = 35.25
$ws.Range("J7").Value = 159.33333
$ws.Range("K7").Value = 35.25
$ws.Range("L7").Value = 159.33333
$ws.Range("M7").Value = 77.75
$ws.Range("N7").Value = -385.33333
$ws.Range("H10").Value = 7143067.5
$ws.Range("I10").Value = 7692457.5
$ws.Range("J10").Value = 1000
$ws.Range("K10").Value = 7692457.5
$ws.Range("L10").Value = 1000
$ws.Range("M10").Value = -7692318.5
$ws.Range("N10").Value = -1278
$ws.Range("H14").Value = 25000
$ws.Range("J14").Value = 25000
$ws.Range("L14").Value = 25000
$ws.Range("N14").Value = -25340
$ws.Range("H15").Value = 2716.7778
$ws.Range("I15").Value = 610.8
$ws.Range("J15").Value = 5349.25
$ws.Range("K15").Value = 610.8
$ws.Range("L15").Value = 5349.25
$ws.Range("M15").Value = -440.8
$ws.Range("N15").Value = -5689.25
$ws.Range("H18").Value = 30000
$ws.Range("J18").Value = 30000
$ws.Range("L18").Value = 30000
$ws.Range("N18").Value = -30460
$ws.Range("H21").Value = 16694.8
$ws.Range("J21").Value = 16694.8
$ws.Range("L21").Value = 16694.8
$ws.Range("N21").Value = -17164.8
$ws.Range("H22").Value = 382.55554
$ws.Range("I22").Value = 185.1
$ws.Range("J22").Value = 629.375
$ws.Range("K22").Value = 185.1
$ws.Range("L22").Value = 629.375
$ws.Range("M22").Value = 164.9
$ws.Range("N22").Value = -1329.375
$ws.Range("H23").Value = 18254.545
$ws.Range("I23").Value = 3000
$ws.Range("J23").Value = 19780
$ws.Range("K23").Value = 3000
$ws.Range("L23").Value = 19780
$ws.Range("M23").Value = -2760
$ws.Range("N23").Value = -20260
$ws.Range("H26").Value = 19500
$ws.Range("J26").Value = 19500
$ws.Range("L26").Value = 19500
$ws.Range("N26").Value = -20074
$ws.Range("H27").Value = 18254.545
$ws.Range("I27").Value = 3000
$ws.Range("J27").Value = 19780
$ws.Range("K27").Value = 3000
$ws.Range("L27").Value = 19780
$ws.Range("M27").Value = -2808
$ws.Range("N27").Value = -20164
$ws.Range("H32").Value = 20499.223
$ws.Range("J32").Value = 21811.625
$ws.Range("L32").Value = 21811.625
$ws.Range("N32").Value = -22443.625
$ws.Range("H33").Value = 11166.125
$ws.Range("I33").Value = 1277
$ws.Range("J33").Value = 17099.6
$ws.Range("K33").Value = 1277
$ws.Range("L33").Value = 17099.6
$ws.Range("M33").Value = -898
$ws.Range("N33").Value = -17857.6
$ws.Range("H36").Value = 11034.417
$ws.Range("I36").Value = 7157
$ws.Range("J36").Value = 14911.833
$ws.Range("K36").Value = 7157
$ws.Range("L36").Value = 14911.833
$ws.Range("M36").Value = -6769
$ws.Range("N36").Value = -15687.833
$ws.Range("H38").Value = 27500
$ws.Range("J38").Value = 27500
$ws.Range("L38").Value = 27500
$ws.Range("N38").Value = -28254
$ws.Range("H39").Value = 15829
$ws.Range("I39").Value = 9697.5
$ws.Range("J39").Value = 19916.666
$ws.Range("K39").Value = 9697.5
$ws.Range("L39").Value = 19916.666
$ws.Range("M39").Value = -9306.5
$ws.Range("N39").Value = -20698.666
$ws.Range("H40").Value = 11034.417
$ws.Range("I40").Value = 7157
$ws.Range("J40").Value = 14911.833
$ws.Range("K40").Value = 7157
$ws.Range("L40").Value = 14911.833
$ws.Range("M40").Value = -6997
$ws.Range("N40").Value = -15231.833
$ws.Range("H44").Value = 37499
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 37499
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 37499
$ws.Range("M44").ClearContents()
$ws.Range("N44").Value = -38383
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("H46").Value = 27500
$ws.Range("J46").Value = 27500
$ws.Range("L46").Value = 27500
$ws.Range("N46").Value = -27922
$ws.Range("H49").Value = 15829
$ws.Range("I49").Value = 9697.5
$ws.Range("J49").Value = 19916.666
$ws.Range("K49").Value = 9697.5
$ws.Range("L49").Value = 19916.666
$ws.Range("M49").Value = -9515.5
$ws.Range("N49").Value = -20280.666
$ws.Range("H50").Value = 13000
$ws.Range("J50").Value = 13000
$ws.Range("L50").Value = 13000
$ws.Range("N50").Value = -14250
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H86").Value = 2321431.8
$ws.Range("I86").Value = 4279073.5
$ws.Range("K86").Value = 4279073.5
$ws.Range("M86").Value = -4277950.5
$ws.Range("H89").Value = 2321431.8
$ws.Range("I89").Value = 4279073.5
$ws.Range("K89").Value = 21395367.5
$ws.Range("M89").Value = -21389751.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 192.16667
$ws.Range("I38").Value = 11
$ws.Range("J38").Value = 228.4
$ws.Range("K38").Value = 33
$ws.Range("L38").Value = 685.2
$ws.Range("M38").Value = 314
$ws.Range("N38").Value = -1379.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5834.816
$ws.Range("I70").Value = 5166.5356
$ws.Range("J70").Value = 7706
$ws.Range("K70").Value = 5166.5356
$ws.Range("L70").Value = 7706
$ws.Range("M70").Value = -4896.5356
$ws.Range("N70").Value = -8246
$ws.Range("H73").Value = 5834.816
$ws.Range("I73").Value = 5166.5356
$ws.Range("J73").Value = 7706
$ws.Range("K73").Value = 5166.5356
$ws.Range("L73").Value = 7706
$ws.Range("M73").Value = -4230.5356
$ws.Range("N73").Value = -9578

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 819.7143
$ws.Range("I46").Value = 948.3333
$ws.Range("J46").Value = 768.26666
$ws.Range("K46").Value = 948.3333
$ws.Range("L46").Value = 768.26666
$ws.Range("M46").Value = -760.3333
$ws.Range("N46").Value = -1144.26666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 6479
$ws.Range("J74").Value = 6479
$ws.Range("L74").Value = 6479
$ws.Range("N74").Value = -8351
$ws.Range("H77").Value = 6479
$ws.Range("J77").Value = 6479
$ws.Range("L77").Value = 19437
$ws.Range("N77").Value = -28797
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("H130").Value = 29476.334
$ws.Range("J130").Value = 29476.334
$ws.Range("L130").Value = 29476.334
